# Revert "Passed" -> "Failed" for the Execution Status column (J3:J5)
# on the TestCases sheet, including the font color (green -> red),
# keeping the existing wrap-text alignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$cells = @("J3", "J4", "J5")
foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.Value2 = "Failed"
    $rng.Font.Color = 255
    $rng.WrapText = $true
}
